$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the data range A1:C78 into an Excel Table ("Table1"), which is
# what Excel does under the hood when you Insert > Table / Ctrl+T over a
# range that already has headers in row 1.
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:C78"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Sort the table by the "Year" column (column C) in descending order,
# moving the newer 2025 scholars to the top of the list while keeping
# 2024 entries (in their original relative order) below.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("C1:C78"), [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlDescending)
$lo.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$lo.Sort.Apply()

# Update the active cell selection to match the saved view state.
$ws.Range("F7").Select()
